$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-16
$data = @{
    2  = @(5, 7)
    3  = @(8, 8)
    4  = @(8, 9)
    5  = @(5, 6)
    6  = @(9, 9)
    7  = @(5, 6)
    8  = @(6, 7)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(6, 6)
    12 = @(7, 7)
    13 = @(5, 6)
    14 = @(4, 6)
    15 = @(5, 5)
    16 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
